$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 2416
$ws.Range("I62").Value = 2249
$ws.Range("J62").Value = 2750
$ws.Range("K62").Value = 2249
$ws.Range("L62").Value = 2750
$ws.Range("M62").Value = -1625
$ws.Range("N62").Value = -3998
$ws.Range("H65").Value = 2416
$ws.Range("I65").Value = 2249
$ws.Range("J65").Value = 2750
$ws.Range("K65").Value = 11245
$ws.Range("L65").Value = 13750
$ws.Range("M65").Value = -8125
$ws.Range("N65").Value = -19990
$ws.Range("H96").Value = 1067
$ws.Range("I96").Value = 970.4
$ws.Range("J96").Value = 1550
$ws.Range("K96").Value = 2911.2
$ws.Range("L96").Value = 4650
$ws.Range("M96").Value = -1538.2
$ws.Range("N96").Value = -7396
$ws.Range("H100").Value = 4751.5
$ws.Range("I100").Value = 4000
$ws.Range("J100").Value = 5002
$ws.Range("K100").Value = 4000
$ws.Range("L100").Value = 5002
$ws.Range("M100").Value = -3459
$ws.Range("N100").Value = -6084
$ws.Range("H107").Value = 707
$ws.Range("I107").Value = 707
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 707
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 1213
$ws.Range("N107").Value = $null
$ws.Range("H132").Value = 7045.4165
$ws.Range("I132").Value = 7354.5
$ws.Range("K132").Value = 22063.5
$ws.Range("M132").Value = -19533.5
$ws.Range("H137").Value = 73277
$ws.Range("I137").Value = 1746.5
$ws.Range("K137").Value = 5239.5
$ws.Range("M137").Value = -2689.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 71432210
$ws.Range("I74").Value = 76926870
$ws.Range("K74").Value = 76926870
$ws.Range("M74").Value = -76925996
$ws.Range("H77").Value = 71432210
$ws.Range("I77").Value = 76926870
$ws.Range("K77").Value = 384634350
$ws.Range("M77").Value = -384629982
$ws.Range("H102").Value = 1813.75
$ws.Range("I102").Value = 1627.5
$ws.Range("K102").Value = 1627.5
$ws.Range("M102").Value = -5.5
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").Value = $null
$ws.Range("H132").Value = 14632.205
$ws.Range("I132").Value = 1717.9166
$ws.Range("J132").Value = 35295.066
$ws.Range("K132").Value = 5153.7498
$ws.Range("L132").Value = 105885.198
$ws.Range("M132").Value = -2623.7498
$ws.Range("N132").Value = -110945.198

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3367.6155
$ws.Range("I20").Value = 4086.7778
$ws.Range("K20").Value = 4086.7778
$ws.Range("M20").Value = -3839.7778
$ws.Range("H70").Value = 105000
$ws.Range("J70").Value = 105000
$ws.Range("L70").Value = 105000
$ws.Range("N70").Value = -105586
$ws.Range("H73").Value = 105000
$ws.Range("J73").Value = 105000
$ws.Range("L73").Value = 105000
$ws.Range("N73").Value = -107028
$ws.Range("H94").Value = 1117.449
$ws.Range("I94").Value = 972.1795
$ws.Range("K94").Value = 972.1795
$ws.Range("M94").Value = -521.1795
$ws.Range("H134").Value = 38831.93
$ws.Range("I134").Value = 50687.547
$ws.Range("K134").Value = 152062.641
$ws.Range("M134").Value = -149527.641

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 25457.215
$ws.Range("I31").Value = 41987.625
$ws.Range("J31").Value = 3416.6667
$ws.Range("K31").Value = 41987.625
$ws.Range("L31").Value = 3416.6667
$ws.Range("M31").Value = -41692.625
$ws.Range("N31").Value = -4006.6667
$ws.Range("H34").Value = 25457.215
$ws.Range("I34").Value = 41987.625
$ws.Range("J34").Value = 3416.6667
$ws.Range("K34").Value = 41987.625
$ws.Range("L34").Value = 3416.6667
$ws.Range("M34").Value = -41785.625
$ws.Range("N34").Value = -3820.6667
$ws.Range("H99").Value = 12198658
$ws.Range("I99").Value = 2876.1538
$ws.Range("J99").Value = 33338012
$ws.Range("K99").Value = 2876.1538
$ws.Range("L99").Value = 33338012
$ws.Range("M99").Value = -1378.1538
$ws.Range("N99").Value = -33341008
$ws.Range("H107").Value = 1070
$ws.Range("I107").Value = 244.42857
$ws.Range("J107").Value = 1895.5714
$ws.Range("K107").Value = 244.42857
$ws.Range("L107").Value = 1895.5714
$ws.Range("M107").Value = 1675.57143
$ws.Range("N107").Value = -5735.5714
$ws.Range("H122").Value = 2692.875
$ws.Range("I122").Value = 3323
$ws.Range("J122").Value = 1642.6666
$ws.Range("K122").Value = 9969
$ws.Range("L122").Value = 4927.9998
$ws.Range("M122").Value = -7519
$ws.Range("N122").Value = -9827.9998
$ws.Range("H126").Value = 12198658
$ws.Range("I126").Value = 2876.1538
$ws.Range("J126").Value = 33338012
$ws.Range("K126").Value = 8628.4614
$ws.Range("L126").Value = 100014036
$ws.Range("M126").Value = -6158.4614
$ws.Range("N126").Value = -100018976
$ws.Range("H134").Value = 1096.1428
$ws.Range("I134").Value = 1031.1666
$ws.Range("J134").Value = 1164.9412
$ws.Range("K134").Value = 3093.4998
$ws.Range("L134").Value = 3494.8236
$ws.Range("M134").Value = -558.4998000000001
$ws.Range("N134").Value = -8564.8236

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 690.375
$ws.Range("J122").Value = 878.6
$ws.Range("L122").Value = 7907.400000000001
$ws.Range("N122").Value = -12807.4
$ws.Range("H129").Value = 358553.44
$ws.Range("J129").Value = 418179
$ws.Range("L129").Value = 1254537
$ws.Range("N129").Value = -1264537
$ws.Range("H131").Value = 778.1900000000001
$ws.Range("J131").Value = 776.9596
$ws.Range("L131").Value = 2330.8788
$ws.Range("N131").Value = -12410.8788

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3979.9714
$ws.Range("I126").Value = 3063.5908
$ws.Range("K126").Value = 9190.7724
$ws.Range("M126").Value = -6720.7724
$ws.Range("H132").Value = 47723.293
$ws.Range("I132").Value = 47588.176
$ws.Range("K132").Value = 142764.528
$ws.Range("M132").Value = -140234.528

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 73.789474
$ws.Range("J55").Value = 113
$ws.Range("L55").Value = 113
$ws.Range("N55").Value = -459
$ws.Range("H69").Value = 40000
$ws.Range("J69").Value = 40000
$ws.Range("L69").Value = 40000
$ws.Range("N69").Value = -41622
$ws.Range("H72").Value = 40000
$ws.Range("J72").Value = 40000
$ws.Range("L72").Value = 120000
$ws.Range("N72").Value = -128112
$ws.Range("H132").Value = 1854.0714
$ws.Range("I132").Value = 1300.7858
$ws.Range("K132").Value = 3902.3574
$ws.Range("M132").Value = -1372.3574
$ws.Range("H136").Value = 39275.848
$ws.Range("I136").Value = 39275.848
$ws.Range("K136").Value = 117827.544
$ws.Range("M136").Value = -115277.544

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 18519962
$ws.Range("I136").Value = 25642266
$ws.Range("J136").Value = 1973.9333
$ws.Range("K136").Value = 76926798
$ws.Range("L136").Value = 5921.7999
$ws.Range("M136").Value = -76924248
